# Sync attendance_reports: fix "Recorded By" ordering so the recorder's
# email appears before "System" (e.g. "System, dnasr281@gmail.com" ->
# "dnasr281@gmail.com, System") in column G of the Session Analysis sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value = $newValue
    }
}
